$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws 'D2' '25.821.48'
Set-TextCell $ws 'D3' '1.706.09'
Set-TextCell $ws 'E3' '  +3.26%  '
Set-TextCell $ws 'E4' '  +0.06%  '
Set-TextCell $ws 'D5' '330.60'
Set-TextCell $ws 'E5' '  +5.62%  '
Set-TextCell $ws 'D6' '0.9986'
Set-TextCell $ws 'E6' '  -0.09%  '
Set-TextCell $ws 'D7' '0.3682'
Set-TextCell $ws 'E7' '  +0.55%  '
Set-TextCell $ws 'D8' '48.48'
Set-TextCell $ws 'E8' '  +3.92%  '
Set-TextCell $ws 'D9' '0.3310'
Set-TextCell $ws 'E9' '  +1.66%  '
Set-TextCell $ws 'D10' '1.170'
Set-TextCell $ws 'E10' '  +3.80%  '
Set-TextCell $ws 'D11' '0.07337'
Set-TextCell $ws 'E11' '  +4.01%  '
Set-TextCell $ws 'D12' '0.9993'
Set-TextCell $ws 'E12' '  +0.16%  '
Set-TextCell $ws 'D13' '6.198'
Set-TextCell $ws 'E13' '  +3.68%  '
Set-TextCell $ws 'D14' '20.00'
Set-TextCell $ws 'E14' '  +2.96%  '
Set-TextCell $ws 'D15' '6.863'
Set-TextCell $ws 'E15' '  +3.60%  '
Set-TextCell $ws 'D16' '1.705.89'
Set-TextCell $ws 'E16' '  +3.04%  '
Set-TextCell $ws 'D17' '0.00001066'
Set-TextCell $ws 'E17' '  +2.07%  '
Set-TextCell $ws 'E18' '  +0.76%  '
Set-TextCell $ws 'D19' '81.04'
Set-TextCell $ws 'E19' '  +2.68%  '
Set-TextCell $ws 'D20' '0.9992'
Set-TextCell $ws 'E20' '  +0.06%  '
Set-TextCell $ws 'B21' 'Uniswap'
Set-TextCell $ws 'C21' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell $ws 'D21' '6.047'
Set-TextCell $ws 'E21' '  +1.80%  '
Set-TextCell $ws 'B22' 'Avalanche'
Set-TextCell $ws 'C22' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell $ws 'D22' '16.16'
Set-TextCell $ws 'E22' '  +2.87%  '
Set-TextCell $ws 'D23' '13.00'
Set-TextCell $ws 'E23' '  +3.23%  '
Set-TextCell $ws 'D24' '25.797.33'
Set-TextCell $ws 'E24' '  +5.54%  '
Set-TextCell $ws 'D25' '2.462'
Set-TextCell $ws 'E25' '  -0.30%  '
Set-TextCell $ws 'D26' '2.474'
Set-TextCell $ws 'E26' '  +6.01%  '
Set-TextCell $ws 'D27' '149.57'
Set-TextCell $ws 'E27' '  +1.84%  '
Set-TextCell $ws 'D28' '19.15'
Set-TextCell $ws 'E28' '  +2.76%  '
Set-TextCell $ws 'D29' '1.284'
Set-TextCell $ws 'E29' '  +6.59%  '
Set-TextCell $ws 'D30' '1.891.71'
Set-TextCell $ws 'E30' '  +2.78%  '
Set-TextCell $ws 'D31' '128.37'
Set-TextCell $ws 'E31' '  +3.22%  '
Set-TextCell $ws 'D32' '4.107'
Set-TextCell $ws 'E32' '  +1.14%  '
Set-TextCell $ws 'D33' '5.933'
Set-TextCell $ws 'E33' '  +3.38%  '
Set-TextCell $ws 'D34' '1.723'
Set-TextCell $ws 'E34' '  +4.13%  '
Set-TextCell $ws 'D35' '0.08518'
Set-TextCell $ws 'E35' '  +0.61%  '
Set-TextCell $ws 'D36' '12.87'
Set-TextCell $ws 'E36' '  +5.74%  '
Set-TextCell $ws 'D37' '5.318'
Set-TextCell $ws 'E37' '  +1.68%  '
Set-TextCell $ws 'D38' '1.276'
Set-TextCell $ws 'E38' '  +0.71%  '
Set-TextCell $ws 'D39' '0.06183'
Set-TextCell $ws 'E39' '  +2.42%  '
Set-TextCell $ws 'D40' '8.511'
Set-TextCell $ws 'E40' '  +4.44%  '
Set-TextCell $ws 'E41' '  +2.08%  '
Set-TextCell $ws 'D42' '0.02249'
Set-TextCell $ws 'E42' '  +0.44%  '
Set-TextCell $ws 'D43' '14.64'
Set-TextCell $ws 'E43' '  +16.58%  '
Set-TextCell $ws 'D44' '0.6105'
Set-TextCell $ws 'E44' '  +3.03%  '
Set-TextCell $ws 'D45' '0.9994'
Set-TextCell $ws 'E45' '  +0.01%  '
Set-TextCell $ws 'D46' '3.839'
Set-TextCell $ws 'E46' '  +1.32%  '
Set-TextCell $ws 'D47' '0.5822'
Set-TextCell $ws 'E47' '  +3.37%  '
Set-TextCell $ws 'D48' '126.77'
Set-TextCell $ws 'E48' '  +2.90%  '
Set-TextCell $ws 'D49' '1.998'
Set-TextCell $ws 'E49' '  +2.30%  '
Set-TextCell $ws 'D50' '0.07221'
Set-TextCell $ws 'E50' '  +4.38%  '
Set-TextCell $ws 'D51' '1.203'
Set-TextCell $ws 'E51' '  +1.37%  '
